$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Mint a fresh (unused) list-numbering definition the same way Word
# does internally whenever list formatting machinery gets touched
# during an editing session -- a throwaway paragraph is numbered and
# then removed again, leaving the newly defined abstract numbering
# behind in numbering.xml without altering the visible document body.
# ------------------------------------------------------------------
$scratch = $d.Paragraphs.Last
$scratch.Range.InsertParagraphAfter()
$scratch = $d.Paragraphs.Last
$scratch.Range.ListFormat.ApplyNumberDefault()
$scratch.Range.Delete()

# ------------------------------------------------------------------
# Fill in the final (previously blank) paragraph with the writeup for
# the new "WoW Character Search" portfolio entry, splitting it into
# four separate paragraphs exactly as in the authored copy.
# ------------------------------------------------------------------
$last = $d.Paragraphs.Last
$last.Range.Text = "I built the front end of a website based on the concept design given to me.`r" + `
    "The website fetches World of Warcraft character data to display it to the user.`r" + `
    "Image showcases the various information boxes that appear when hovering over certain elements.`r" + `
    "Concept Design Image, website was built using JS/HTML/CSS."
